# Updates cryptos list prices/volumes (rows 2-50) and replaces EOS with Aave (row 51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.637.42"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "1.854.60"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'264.90"
$ws.Range("E5").Value = "  +2.18%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "'0.5234"
$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("D8").Value = "'0.3287"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.06801"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").Value = "'18.87"
$ws.Range("E10").Value = "  -3.06%  "

$ws.Range("D11").Value = "'0.7774"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").Value = "'0.07726"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").Value = "1.850.29"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "'88.66"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").Value = "'5.039"
$ws.Range("E15").Value = "  -0.54%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "'0.000007991"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").Value = "26.664.52"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "2.086.85"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").Value = "'9.552"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").Value = "'6.008"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D26").Value = "'2.200"
$ws.Range("E26").Value = "  -5.89%  "

$ws.Range("D27").Value = "'1.680"
$ws.Range("E27").Value = "  +2.17%  "

$ws.Range("D28").Value = "'16.99"
$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("D29").Value = "'112.41"
$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("D30").Value = "'4.194"
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("E31").Value = "  -1.04%  "

$ws.Range("D32").Value = "'0.08762"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "'0.04836"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("D35").Value = "'2.868"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").Value = "'0.7150"
$ws.Range("E36").Value = "  +0.73%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "'0.01786"
$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").Value = "'0.4890"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("D41").Value = "'113.00"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").Value = "'0.9014"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").Value = "'6.081"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").Value = "'7.734"
$ws.Range("E45").Value = "  -1.07%  "

$ws.Range("D46").Value = "'0.4206"
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("D47").Value = "'9.126"
$ws.Range("E47").Value = "  -1.25%  "

$ws.Range("D48").Value = "'0.05924"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").Value = "'0.1241"
$ws.Range("E49").Value = "  -4.13%  "

$ws.Range("D50").Value = "'35.08"
$ws.Range("E50").Value = "  -0.78%  "

# Row 51: coin swapped from EOS to Aave (name, link, price, volume)
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'60.13"
$ws.Range("E51").Value = "  +0.67%  "
